$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-14 from 2023-09-13 (45182)
# to 2023-09-15 (45184), matching the new date serial value.
$ws.Range("C2:C14").Value = 45184
